# "Testing with new leds"
# Update the calibration sheet (20250309) with new LED test data in column G,
# rename the G1 test-run label, and remove a leftover duplicate chart anchor.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("20250309")

# G1: rename the test-run identifier.
$ws.Range("G1").Value = "595xDM605xF600"

# G5:G20 - new photometer readings recorded for this LED run.
$ws.Range("G5").Value  = 0
$ws.Range("G6").Value  = 0.069
$ws.Range("G7").Value  = 0.172
$ws.Range("G8").Value  = 0.33
$ws.Range("G9").Value  = 1.42
$ws.Range("G10").Value = 2.56
$ws.Range("G11").Value = 3.72
$ws.Range("G12").Value = 5.42
$ws.Range("G13").Value = 8.21
$ws.Range("G14").Value = 10.9
$ws.Range("G15").Value = 13.47
$ws.Range("G16").Value = 15.9
$ws.Range("G17").Value = 18.14
$ws.Range("G18").Value = 20.22
$ws.Range("G19").Value = 21.98
$ws.Range("G20").Value = 23.43

# Row 22 header: column G now reports against the MEA3_625 channel.
$ws.Range("G22").Value = "MEA3_625"

# Row 23: direct measurement for the new channel.
$ws.Range("G23").Value = 7500

# Row 25 had no data for this run - clear the stray formatted blank cell.
$ws.Range("G25").Clear()

# Remove the leftover duplicate chart placeholder (no graphic, same block as
# the "Graphique 5_0" chart already placed to the left).
$shapes = $ws.Shapes
for ($i = $shapes.Count; $i -ge 1; $i--) {
    $shape = $shapes.Item($i)
    if ($shape.Name -eq "") {
        $shape.Delete()
    }
}
